# Apply the "revert admin dev default" edit:
#  - Weekly Timesheet: rename clients in B2:B6, zero out Rate (E) and Total (F)
#    for the detail rows and the subtotal/hourly-subtotal/grand-total rows.
#  - Jason Schema: same client renames (C/D stay as Date; client is col D),
#    zero out Rate (F) / Total (G), and clear the seeded Notes column (I).
#  - Employee ID changes from the old seeded GUID-style id to a short id.

$wb = $excel.ActiveWorkbook

$timesheet = $wb.Worksheets.Item("Weekly Timesheet")
$schema = $wb.Worksheets.Item("Jason Schema")

# New client names, in row order (rows 2..6)
$newClients = @("Evans", "Oglesby", "Muncey", "Lucas", "Bailey")

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $name = $newClients[$i]

    # Weekly Timesheet sheet: B = Client, E = Rate, F = Total
    $timesheet.Range("B$row").Value = $name
    $timesheet.Range("E$row").Value = 0
    $timesheet.Range("F$row").Value = 0

    # Jason Schema sheet: B = Employee ID, D = Client, F = Rate, G = Total, I = Notes
    $schema.Range("B$row").Value = "emp_pw6be4hd"
    $schema.Range("D$row").Value = $name
    $schema.Range("F$row").Value = 0
    $schema.Range("G$row").Value = 0

    # Clear the seeded "Notes" text but leave a (typed) empty string behind,
    # same shape as the shared empty-string cell used elsewhere in the sheet.
    $schema.Range("I$row").Value = "'"
    $schema.Range("I$row").Style = "Normal"
}

# Weekly Timesheet subtotal / hourly-subtotal / grand-total Total cells
$timesheet.Range("F8").Value = 0
$timesheet.Range("F11").Value = 0
$timesheet.Range("F13").Value = 0
